# Fix the sheet name: "自作関数" -> "test_data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "test_data"
